# Update market price / profit figures on the Goblin Profits sheets
# (refreshed data from the scheduled market-data runner)
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 166667570
$ws.Range("I29").Value = 500000000
$ws.Range("J29").Value = 1360.5
$ws.Range("K29").Value = 1500000000
$ws.Range("L29").Value = 4081.5
$ws.Range("M29").Value = -1499999719
$ws.Range("N29").Value = -4643.5
# Row 43
$ws.Range("H43").Value = 6624.875
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 9399.799999999999
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 9399.799999999999
$ws.Range("M43").Value = -1931
$ws.Range("N43").Value = -9537.799999999999
# Row 58
$ws.Range("H58").Value = 14723984
$ws.Range("I58").Value = 50000150
$ws.Range("J58").Value = 25579.916
$ws.Range("K58").Value = 150000450
$ws.Range("L58").Value = 76739.74800000001
$ws.Range("M58").Value = -150000300
$ws.Range("N58").Value = -77039.74800000001
# Row 137
$ws.Range("H137").Value = 9772.462
$ws.Range("I137").Value = 12364.2
$ws.Range("J137").Value = 1133.3334
$ws.Range("K137").Value = 37092.60000000001
$ws.Range("L137").Value = 3400.0002
$ws.Range("M137").Value = -34542.60000000001
$ws.Range("N137").Value = -8500.0002
# Row 138
$ws.Range("H138").Value = 6255.4873
$ws.Range("I138").Value = 3099
$ws.Range("J138").Value = 7070.0645
$ws.Range("K138").Value = 9297
$ws.Range("L138").Value = 21210.1935
$ws.Range("M138").Value = -4157
$ws.Range("N138").Value = -31490.1935

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1598.56
$ws.Range("I2").Value = 572.0833
$ws.Range("J2").Value = 2546.077
$ws.Range("K2").Value = 572.0833
$ws.Range("L2").Value = 2546.077
$ws.Range("M2").Value = -459.0833
$ws.Range("N2").Value = -2772.077
# Row 32
$ws.Range("H32").Value = 2367.9092
$ws.Range("I32").Value = 2374.7222
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 2374.7222
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -2087.7222
$ws.Range("N32").Value = -2574
# Row 61
$ws.Range("H61").Value = 2552.578
$ws.Range("I61").Value = 2544.95
$ws.Range("J61").Value = 2613.6
$ws.Range("K61").Value = 2544.95
$ws.Range("L61").Value = 2613.6
$ws.Range("M61").Value = -2332.95
$ws.Range("N61").Value = -3037.6
# Row 116
$ws.Range("H116").Value = 1598.56
$ws.Range("I116").Value = 572.0833
$ws.Range("J116").Value = 2546.077
$ws.Range("K116").Value = 572.0833
$ws.Range("L116").Value = 2546.077
$ws.Range("M116").Value = 1721.9167
$ws.Range("N116").Value = -7134.077
# Row 132
$ws.Range("H132").Value = 2675.5
$ws.Range("I132").Value = 2404.4285
$ws.Range("J132").Value = 3149.875
$ws.Range("K132").Value = 7213.2855
$ws.Range("L132").Value = 9449.625
$ws.Range("M132").Value = -4683.2855
$ws.Range("N132").Value = -14509.625
# Row 136
$ws.Range("H136").Value = 2552.578
$ws.Range("I136").Value = 2544.95
$ws.Range("J136").Value = 2613.6
$ws.Range("K136").Value = 7634.849999999999
$ws.Range("L136").Value = 7840.799999999999
$ws.Range("M136").Value = -5084.849999999999
$ws.Range("N136").Value = -12940.8

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1598.56
$ws.Range("I3").Value = 572.0833
$ws.Range("J3").Value = 2546.077
$ws.Range("K3").Value = 572.0833
$ws.Range("L3").Value = 2546.077
$ws.Range("M3").Value = -458.0833
$ws.Range("N3").Value = -2774.077
# Row 88
$ws.Range("H88").Value = 16060.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 16060.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 16060.5
$ws.Range("N88").Value = -16872.5
# Row 91
$ws.Range("H91").Value = 16060.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 16060.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 16060.5
$ws.Range("N91").Value = -18868.5
# Row 94
$ws.Range("H94").Value = 87298.14999999999
$ws.Range("I94").Value = 105683.03
$ws.Range("J94").Value = 626.5714
$ws.Range("K94").Value = 105683.03
$ws.Range("L94").Value = 626.5714
$ws.Range("M94").Value = -105232.03
$ws.Range("N94").Value = -1528.5714
# Row 99
$ws.Range("H99").Value = 4960.9165
$ws.Range("I99").Value = 4342.4443
$ws.Range("J99").Value = 6816.3335
$ws.Range("K99").Value = 4342.4443
$ws.Range("L99").Value = 6816.3335
$ws.Range("M99").Value = -2844.4443
$ws.Range("N99").Value = -9812.333500000001
# Row 107
$ws.Range("H107").Value = 8800.454
$ws.Range("I107").Value = 7654.5386
$ws.Range("J107").Value = 9545.299999999999
$ws.Range("K107").Value = 7654.5386
$ws.Range("L107").Value = 9545.299999999999
$ws.Range("M107").Value = -5734.5386
$ws.Range("N107").Value = -13385.3
# Row 134
$ws.Range("H134").Value = 2132.4092
$ws.Range("I134").Value = 2070.8948
$ws.Range("J134").Value = 2522
$ws.Range("K134").Value = 6212.6844
$ws.Range("L134").Value = 7566
$ws.Range("M134").Value = -3677.6844
$ws.Range("N134").Value = -12636

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 3272.2856
$ws.Range("I122").Value = 3130.7144
$ws.Range("J122").Value = 3413.8572
$ws.Range("K122").Value = 9392.143199999999
$ws.Range("L122").Value = 10241.5716
$ws.Range("M122").Value = -6942.143199999999
$ws.Range("N122").Value = -15141.5716

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1200
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1031
$ws.Range("N22").ClearContents()
# Row 27
$ws.Range("H27").Value = 400
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1200
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -1098
$ws.Range("N27").ClearContents()
# Row 59
$ws.Range("H59").Value = 4000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 4000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 12000
$ws.Range("N59").Value = -13080
# Row 92
$ws.Range("H92").Value = 1764.6
$ws.Range("I92").Value = 2050.125
$ws.Range("J92").Value = 622.5
$ws.Range("K92").Value = 6150.375
$ws.Range("L92").Value = 1867.5
$ws.Range("M92").Value = -4902.375
$ws.Range("N92").Value = -4363.5
# Row 97
$ws.Range("H97").Value = 217.6842
$ws.Range("I97").Value = 190
$ws.Range("J97").Value = 222.875
$ws.Range("K97").Value = 570
$ws.Range("L97").Value = 668.625
$ws.Range("M97").Value = -74
$ws.Range("N97").Value = -1660.625
# Row 113
$ws.Range("H113").Value = 1521.2858
$ws.Range("I113").Value = 537.5
$ws.Range("J113").Value = 2833
$ws.Range("K113").Value = 1612.5
$ws.Range("L113").Value = 8499
$ws.Range("M113").Value = 557.5
$ws.Range("N113").Value = -12839

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1462.7858
$ws.Range("I102").Value = 1287.7778
$ws.Range("J102").Value = 1777.8
$ws.Range("K102").Value = 1287.7778
$ws.Range("L102").Value = 1777.8
$ws.Range("M102").Value = 334.2221999999999
$ws.Range("N102").Value = -5021.8
# Row 132
$ws.Range("H132").Value = 2982.3
$ws.Range("I132").Value = 3149.25
$ws.Range("J132").Value = 2731.875
$ws.Range("K132").Value = 9447.75
$ws.Range("L132").Value = 8195.625
$ws.Range("M132").Value = -6917.75
$ws.Range("N132").Value = -13255.625

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 19
$ws.Range("H19").Value = 2027.8
$ws.Range("I19").Value = 2284.75
$ws.Range("J19").Value = 1000
$ws.Range("K19").Value = 2284.75
$ws.Range("L19").Value = 1000
$ws.Range("M19").Value = -2114.75
$ws.Range("N19").Value = -1340
# Row 68
$ws.Range("H68").Value = 4769.8423
$ws.Range("I68").Value = 2901.9285
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 2901.9285
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = -2152.9285
$ws.Range("N68").Value = -11498
# Row 71
$ws.Range("H71").Value = 4769.8423
$ws.Range("I71").Value = 2901.9285
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 14509.6425
$ws.Range("L71").Value = 50000
$ws.Range("M71").Value = -10765.6425
$ws.Range("N71").Value = -57488
# Row 122
$ws.Range("H122").Value = 4320.4116
$ws.Range("I122").Value = 4376.6
$ws.Range("J122").Value = 3899
$ws.Range("K122").Value = 13129.8
$ws.Range("L122").Value = 11697
$ws.Range("M122").Value = -10679.8
$ws.Range("N122").Value = -16597
# Row 133
$ws.Range("H133").Value = 88750
$ws.Range("I133").Value = 130000
$ws.Range("J133").Value = 75000
$ws.Range("K133").Value = 130000
$ws.Range("L133").Value = 75000
$ws.Range("M133").Value = -127470
$ws.Range("N133").Value = -80060

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1861.2222
$ws.Range("I81").Value = 1740.2
$ws.Range("J81").Value = 2012.5
$ws.Range("K81").Value = 3480.4
$ws.Range("L81").Value = 4025
$ws.Range("M81").Value = -2419.4
$ws.Range("N81").Value = -6147
# Row 84
$ws.Range("H84").Value = 1861.2222
$ws.Range("I84").Value = 1740.2
$ws.Range("J84").Value = 2012.5
$ws.Range("K84").Value = 17402
$ws.Range("L84").Value = 20125
$ws.Range("M84").Value = -12098
$ws.Range("N84").Value = -30733
# Row 113
$ws.Range("H113").Value = 1898.5714
$ws.Range("I113").Value = 1218.8889
$ws.Range("J113").Value = 3122
$ws.Range("K113").Value = 3656.6667
$ws.Range("L113").Value = 9366
$ws.Range("M113").Value = -1486.6667
$ws.Range("N113").Value = -13706
# Row 122
$ws.Range("H122").Value = 3770.75
$ws.Range("I122").Value = 4178.3
$ws.Range("J122").Value = 1733
$ws.Range("K122").Value = 12534.9
$ws.Range("L122").Value = 5199
$ws.Range("M122").Value = -10084.9
$ws.Range("N122").Value = -10099
